# Add a new course entry to the "Tabel1" table on the "index" sheet:
# "Introductory teacher training" at Ghent University, 2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item("Tabel1")
$newRow = $tbl.ListRows.Add()

$newRow.Range.Item(1, 1).Value = "Introductory teacher training"
$newRow.Range.Item(1, 2).Value = "Ghent University"
$newRow.Range.Item(1, 3).Value = 2019

$ws.Range("A17").Select()
